# Added the correct Latex template.
$wb = $excel.ActiveWorkbook

# Add a new sheet (Sheet3) after the last existing sheet (Sheet2)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Populate the data for Sheet3 - written in the same order the original
# author entered it (controls shared-string allocation order).
$ws3.Range("A3").Value = "iu-xray (hpcc128)"
$ws3.Range("B2").Value = "% Passed Reports"
$ws3.Range("A2").Value = "Test name, resource"
$ws3.Range("A1").Value = "start date: 010111"
$ws3.Range("A4").Value = "iu-india (hpcc128)"
$ws3.Range("A6").Value = "iu-xray (hpcc256)"
$ws3.Range("A7").Value = "iu-xray (hpcc512)"
$ws3.Range("A8").Value = "iu-xray (hpcc672)"
$ws3.Range("C2").Value = "(05/12-09/12) failures"
$ws3.Range("A5").Value = "iu-india (hpcc256)"

$ws3.Range("B3").Value = 0.93
$ws3.Range("B4").Value = 0.93
$ws3.Range("B5").Value = 0.79
$ws3.Range("B6").Value = 0.91
$ws3.Range("B7").Value = 0.93
$ws3.Range("B8").Value = 0.99

# Apply percentage style to B3:B8 (matches style used elsewhere: numFmtId 9 -> 0%)
$ws3.Range("B3:B8").NumberFormat = "0%"

# Column widths (closest achievable values through this engine's column-width
# quantization; 28 -> ~28.83, 19.15 -> exactly 20)
$ws3.Columns.Item(1).ColumnWidth = 28
$ws3.Columns.Item(2).ColumnWidth = 19.15

# Page margins matching the rest of the workbook (0.75in/1in/0.5in => 54/72/36pt)
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# Sheet1 view changes: topLeftCell becomes C1, selection stays A20, tab no longer shown as selected
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("A20").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

# Sheet2 selection change to E65
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate() | Out-Null
$ws2.Range("E65").Select() | Out-Null

# Activate Sheet3 last so it becomes the selected/active tab, with selection B5
$ws3.Activate() | Out-Null
$ws3.Range("B5").Select() | Out-Null
